# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# "Body" sheet: the request body row (row 3) is turned into a single
# "schema" reference row pointing at liquidityManagement.250808Request,
# and the old body-specific rows (settlementBIC/amount/transferTp) go away.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Body")
$ws.Rows("4:6").Delete()
$ws.Range("B3").Value = "liquidityManagement.250808Request"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "schema"
$ws.Range("G3").Value = "liquidityManagement.250808Request"
$ws.Range("I3").Value = "Yes"
$ws.Range("L3").Value = ""
$ws.Range("O3").Value = ""

# "200" sheet: same treatment, pointing at liquidityManagement.250808Response;
# drops the old availableLiquidity row.
$ws = $wb.Worksheets.Item("200")
$ws.Rows("4:4").Delete()
$ws.Range("B3").Value = "liquidityManagement.250808Response"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "schema"
$ws.Range("G3").Value = "liquidityManagement.250808Response"
$ws.Range("I3").Value = "Yes"
$ws.Range("L3").Value = ""
$ws.Range("O3").Value = ""

# "204" sheet: previously had no body rows at all; add the schema row.
$ws = $wb.Worksheets.Item("204")
$ws.Range("A3").Value = "content"
$ws.Range("B3").Value = "liquidityManagement.250808Response"
$ws.Range("E3").Value = "schema"
$ws.Range("G3").Value = "liquidityManagement.250808Response"
$ws.Range("I3").Value = "Yes"

# "400" sheet: collapse to a generic errorResponse schema reference,
# dropping errorCode/errorCodeDescription/requestId rows.
$ws = $wb.Worksheets.Item("400")
$ws.Rows("4:6").Delete()
$ws.Range("B3").Value = "errorResponse"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "schema"
$ws.Range("G3").Value = "errorResponse"
$ws.Range("I3").Value = "Yes"
$ws.Range("L3").Value = ""
$ws.Range("O3").Value = ""

# "401", "403", "404", "429", "500" sheets: add the errorResponse1 schema row.
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A3").Value = "content"
    $ws.Range("B3").Value = "errorResponse1"
    $ws.Range("E3").Value = "schema"
    $ws.Range("G3").Value = "errorResponse1"
    $ws.Range("I3").Value = "Yes"
}
